$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide "Stock assessments with TMB": merge the two short paragraphs
#    "So can we do assessments in TMB?" / "Of course! " into a single
#    paragraph "So can we do assessments in TMB? Yes! ".
# ---------------------------------------------------------------------------
$cr = [char]13
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            $full = $tr.Text
            $old = "So can we do assessments in TMB?" + $cr + "Of course! " + $cr
            if ($full.Contains($old)) {
                $new = "So can we do assessments in TMB? Yes! " + $cr
                $tr.Text = $full.Replace($old, $new)
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide "Why TMB instead of ADMB?": fix the fixed hyper-variance example
#    from "=.5) instead of estimated" to "=.1) instead of estimated".
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            $full = $tr.Text
            $target = "=.5) instead of estimated"
            $idx0 = $full.IndexOf($target)
            if ($idx0 -ge 0) {
                $c = $tr.Characters($idx0 + 1, $target.Length)
                $c.Text = "=.1) instead of estimated"
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 3) Slide "Recap": the content placeholder grew slightly taller and moved
#    up (autofit re-flow), keeping Left/Width fixed.
#       old: off x=628650  y=1417638   ext cx=7886700 cy=4667249
#       new: off x=628650  y=1201272   ext cx=7886700 cy=4883616
#    Shape.Top/Height are EMU/12700 (points) on a lossy single-precision
#    float -- nudge by a hair so the round-trip lands on the exact EMU.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.HasTextFrame) {
            $t = $sh.TextFrame.TextRange.Text
            if ($t -like "*Moving beyond GLMMs shows the real power of TMB*") {
                $sh.Top = (1201272 / 12700.0) + 0.00001
                $sh.Height = (4883616 / 12700.0) + 0.00001
            }
        }
    }
}
